# Update crypto price (column D) and volume-change (column E) values
# Values are written with a leading apostrophe to force Excel to store them
# as literal text (preserving exact formatting such as trailing zeros and
# the dotted "24.415.02" style price strings) rather than re-parsing them as
# numbers. The Style reset clears the transient "quote prefix" text format
# that gets attached when the text-like-a-number heuristic kicks in, so the
# cell ends up on the same default style as before the edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.459.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.59%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.650.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -3.52%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.43%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'311.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.34%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.9990"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.04%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.3651"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -2.89%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'46.64"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -6.16%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.3249"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -6.28%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'1.126"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -7.51%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07036"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -7.29%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.9980"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.40%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'5.971"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -6.02%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'19.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -9.10%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'6.607"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -6.82%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'1.650.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -3.72%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.00001041"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -8.32%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.06576"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.60%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.9986"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.02%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'78.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -7.58%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'5.933"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -7.59%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'15.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -10.00%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'12.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -5.54%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'24.426.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.90%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.455"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.32%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.328"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -17.02%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'146.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.99%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'18.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -9.61%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.831.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.80%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.189"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -5.00%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'123.75"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -7.21%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'4.085"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.47%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'5.701"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -17.68%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.08438"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -4.93%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.652"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -6.23%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'12.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -13.03%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'5.219"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -8.20%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.261"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.79%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.06024"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -9.95%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.02229"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -7.90%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.2061"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -8.20%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'8.155"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -13.01%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.9983"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.11%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.5904"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -9.01%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'3.763"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -2.02%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'12.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -9.93%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.5621"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -9.27%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'123.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -5.89%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.950"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -9.07%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.06902"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'1.182"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -3.57%  "
$ws.Range("E51").Style = "Normal"
